$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K13").Value = -0.0663
$ws.Range("L13").Value = -0.0079
$ws.Range("M13").Value = -0.0126
$ws.Range("N13").Value = -0.0183
$ws.Range("O13").Value = -0.0516
$ws.Range("P13").Value = 0.0033
$ws.Range("Q13").Value = 0.0125
$ws.Range("R13").Value = 0.012
$ws.Range("S13").Value = 0.0213
$ws.Range("I16").Value = -0.0637
$ws.Range("J17").Value = -0.0648
$ws.Range("K17").Value = -0.0475
$ws.Range("L17").Value = -0.1033
$ws.Range("M17").Value = -0.0365
$ws.Range("N17").Value = -0.0239
$ws.Range("O17").Value = -0.0264
$ws.Range("P17").Value = -0.0206
$ws.Range("Q17").Value = -0.0218
$ws.Range("R17").Value = -0.0304
$ws.Range("S17").Value = -0.0213
$ws.Range("I30").Value = -2.486
$ws.Range("J31").Value = -1.0237
$ws.Range("K31").Value = -1.4718
$ws.Range("L31").Value = -1.7899
$ws.Range("M31").Value = -0.8836
$ws.Range("N31").Value = -0.7029
$ws.Range("O31").Value = -0.7425
$ws.Range("P31").Value = -0.3441
$ws.Range("Q31").Value = -0.3357
$ws.Range("R31").Value = -0.0937
$ws.Range("S31").Value = 0.0342
$ws.Range("L39").Value = -0.4285
$ws.Range("M39").Value = -0.0125
$ws.Range("N39").Value = -0.0123
$ws.Range("O39").Value = -0.6088
$ws.Range("P39").Value = -0.0001
$ws.Range("Q39").Value = -0.0001
$ws.Range("R39").Value = -0.0062
$ws.Range("K47").Value = -0.0025
$ws.Range("L47").Value = 0.004
$ws.Range("M47").Value = 0.12
$ws.Range("N47").Value = 0.1118
$ws.Range("O47").Value = 0.0948
$ws.Range("P47").Value = -0.01
$ws.Range("Q47").Value = -0.0672
$ws.Range("R47").Value = -0.055
$ws.Range("S47").Value = -0.0536
$ws.Range("K69").Value = 0.0176
$ws.Range("L69").Value = -0.0715
$ws.Range("M69").Value = 0.0232
$ws.Range("N69").Value = 0.003
$ws.Range("P69").Value = 0.0739
$ws.Range("Q69").Value = -0.0224
$ws.Range("R69").Value = -0.0214
$ws.Range("S69").Value = -0.0205
$ws.Range("I72").Value = 0.0765
$ws.Range("J73").Value = 0.0259
$ws.Range("K73").Value = 0.0149
$ws.Range("L73").Value = -0.0646
$ws.Range("M73").Value = -0.0207
$ws.Range("N73").Value = -0.0061
$ws.Range("O73").Value = 0.013
$ws.Range("P73").Value = 0.0137
$ws.Range("Q73").Value = 0.0066
$ws.Range("R73").Value = -0.0115
$ws.Range("S73").Value = -0.0178
$ws.Range("I86").Value = 0.0769
$ws.Range("J87").Value = 1.078
$ws.Range("K87").Value = 0.3607
$ws.Range("L87").Value = 0.7619
$ws.Range("M87").Value = 0.1807
$ws.Range("N87").Value = -0.2328
$ws.Range("O87").Value = -0.5743
$ws.Range("P87").Value = 0.024
$ws.Range("Q87").Value = 0.0026
$ws.Range("R87").Value = -0.0864
$ws.Range("S87").Value = -0.0811
$ws.Range("L95").Value = 0.6121
$ws.Range("M95").Value = -0.0122
$ws.Range("N95").Value = -0.0121
$ws.Range("O95").Value = -0.5982
$ws.Range("P95").Value = -0.0001
$ws.Range("Q95").Value = -0.0001
$ws.Range("R95").Value = -0.0062
$ws.Range("K103").Value = -0.0074
$ws.Range("L103").Value = 0.0856
$ws.Range("M103").Value = -0.0053
$ws.Range("N103").Value = 0.0009
$ws.Range("O103").Value = -0.002
$ws.Range("P103").Value = -0.0929
$ws.Range("Q103").Value = -0.0074
$ws.Range("R103").Value = -0.0072
$ws.Range("S103").Value = -0.0069
